$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

# Update stock quantities
$ws.Range("D7").Value = 9850
$ws.Range("D9").Value = 100
$ws.Range("D11").Value = 6

# Update selected cell to reflect where the user left off
$ws.Activate()
$ws.Range("G11").Select()
